# Duplicated catalogue values for NatRel deleted
#
# The NatRelItem sheet contained two copies of the same 9 catalogue rows
# (ids 73-81 duplicating ids 122-130, sharing the same code/name pairs but
# different GUIDs). Remove the older, superseded duplicates and keep the
# later set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NatRelItem")

# Old row numbers (1-indexed, header is row 1) that are plain duplicates of
# rows further down the sheet and must be removed. Row 8 (T_Id 79,
# "geolWasteDisp") is the one duplicate that is kept, so it is not listed.
$rowsToDelete = @(2, 3, 4, 5, 6, 7, 9, 10, 17)

# Delete from the bottom up so earlier row numbers stay valid while we work.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

# Restore the header-aware selection Excel leaves on the sheet, then move
# the active tab back to the first worksheet (AssetKindItem).
$ws.Activate()
$ws.Range("D29").Select()

$wb.Worksheets.Item(1).Activate()
